$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D values are stored as text (e.g. thousand-separator-by-dot
# style numbers, or decimals that must keep trailing zeros), so we force
# text type on write to avoid Excel auto-converting them to numeric values,
# then restore the default "Normal" style so no stray style index remains
# attached to the cell.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.753.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.34%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.304.29'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.06%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.34%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.502'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.31%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.05%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.85%  '

# Row 11
$ws.Range("E11").Value = '  +0.37%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.49%  '

# Row 13
$ws.Range("E13").Value = '  +2.28%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.30%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.659.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.20%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.297.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.88%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.800'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.73%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.655.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.32%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.87%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0900'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.30%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.40%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.16%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.62'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.41%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.08%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.46%  '

# Row 27
$ws.Range("E27").Value = '  -0.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.42%  '

# Row 29
$ws.Range("E29").Value = '  +0.35%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.59'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.91%  '

# Row 31
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.62%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.05%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.05%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.77'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.93%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.34%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.25'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.61%  '

# Row 37
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.52%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0699'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.46%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.44%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1000'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.07%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.86%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.110'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.20%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.05%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.967.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.26%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0281'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.54%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.57%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.01%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.60%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.524.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.22%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.65%  '

# Row 51
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.01%  '
